# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoffs have now been handed back (translated files returned and
# back in sync with en-US):
#   - Status changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every tracked source file.
#   - The "Latest Handback DateTime" placeholder (0001-01-01 00:00:00) is
#     replaced with the real timestamp the handback was received.
#   - Two new columns of data are populated on each language sheet:
#       F "Latest Target File"   - hyperlink back to the source file
#       G "Latest Handback File" - hyperlink to the returned/translated file
#
$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
# The Overview sheet's zh-cn/de-de status cells share the same underlying
# string as the per-language sheets, so refresh them too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("H2").Value = "2016-03-22 14:50:05"
$wsZh.Range("H3").Value = "2016-03-22 14:50:05"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0bdc3ce0770b96609637505421ea248c84d88b87/e2e/01f0c69f-aee0-4b1c-811a-005efcd9f20c.md",
    "",
    "",
    "01f0c69f-aee0-4b1c-811a-005efcd9f20c.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5f92b229c8e1df73d2dbbc2cd7f732588f1135d1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/01f0c69f-aee0-4b1c-811a-005efcd9f20c.fd2ed005646ee695aff98991dd2c08ae392f2cca.zh-cn.xlf",
    "",
    "",
    "01f0c69f-aee0-4b1c-811a-005efcd9f20c.fd2ed005646ee695aff98991dd2c08ae392f2cca.zh-cn.xlf"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0bdc3ce0770b96609637505421ea248c84d88b87/e2e/d41011ff-72cf-4953-909e-a023866d6408.md",
    "",
    "",
    "d41011ff-72cf-4953-909e-a023866d6408.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5f92b229c8e1df73d2dbbc2cd7f732588f1135d1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/d41011ff-72cf-4953-909e-a023866d6408.cb1834f66b68cff4dfa34822e2edf4d0629e97dd.zh-cn.xlf",
    "",
    "",
    "d41011ff-72cf-4953-909e-a023866d6408.cb1834f66b68cff4dfa34822e2edf4d0629e97dd.zh-cn.xlf"
) | Out-Null

# --- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("H2").Value = "2016-03-22 14:50:19"
$wsDe.Range("H3").Value = "2016-03-22 14:50:19"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0bdc3ce0770b96609637505421ea248c84d88b87/e2e/01f0c69f-aee0-4b1c-811a-005efcd9f20c.md",
    "",
    "",
    "01f0c69f-aee0-4b1c-811a-005efcd9f20c.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec7dd2e803d4dfd525b7b0a642a3161b7e6003a8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/01f0c69f-aee0-4b1c-811a-005efcd9f20c.fd2ed005646ee695aff98991dd2c08ae392f2cca.de-de.xlf",
    "",
    "",
    "01f0c69f-aee0-4b1c-811a-005efcd9f20c.fd2ed005646ee695aff98991dd2c08ae392f2cca.de-de.xlf"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0bdc3ce0770b96609637505421ea248c84d88b87/e2e/d41011ff-72cf-4953-909e-a023866d6408.md",
    "",
    "",
    "d41011ff-72cf-4953-909e-a023866d6408.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec7dd2e803d4dfd525b7b0a642a3161b7e6003a8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/d41011ff-72cf-4953-909e-a023866d6408.cb1834f66b68cff4dfa34822e2edf4d0629e97dd.de-de.xlf",
    "",
    "",
    "d41011ff-72cf-4953-909e-a023866d6408.cb1834f66b68cff4dfa34822e2edf4d0629e97dd.de-de.xlf"
) | Out-Null
